$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before the existing row 902, shifting rows 902:988 down to 907:993.
$ws.Rows("902:906").Insert()

# New weekly price rows (Tomate, Comercializadora del Agro de Limari, Coquimbo),
# mirroring the shared/common columns already used throughout the sheet.
$mercado = "Comercializadora del Agro de Limarí"
$region = "Coquimbo"
$categoriaId = 100112020
$categoria = "Tomate"
$unidad = "`$/bandeja 18 kilos"
$origen = "Provincia de Limarí"
$kgUnidades = 18
$clasificacion = "Hortaliza"

function Set-PrecioRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [double]$PrecioKg
    )

    $ws.Cells.Item($Row, 1).Value = 2
    $ws.Cells.Item($Row, 2).Value = $mercado
    $ws.Cells.Item($Row, 3).Value = $region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 4
    $ws.Cells.Item($Row, 6).Value = $categoriaId
    $ws.Cells.Item($Row, 7).Value = $categoria
    $ws.Cells.Item($Row, 8).Value = $Variedad
    $ws.Cells.Item($Row, 9).Value = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $unidad
    $ws.Cells.Item($Row, 15).Value = $origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $kgUnidades
    $ws.Cells.Item($Row, 18).Value = $clasificacion
}

Set-PrecioRow -Row 902 -Fecha 45106 -Variedad "Larga vida" -Calidad "Primera"  -Volumen 2200 -PrecioMin 10000 -PrecioMax 11000 -PrecioProm 10500 -PrecioKg 583
Set-PrecioRow -Row 903 -Fecha 45106 -Variedad "Larga vida" -Calidad "Segunda"  -Volumen 1800 -PrecioMin 8000  -PrecioMax 9000  -PrecioProm 8500  -PrecioKg 472
Set-PrecioRow -Row 904 -Fecha 45106 -Variedad "Larga vida" -Calidad "Tercera"  -Volumen 700  -PrecioMin 5000  -PrecioMax 6000  -PrecioProm 5500  -PrecioKg 306
Set-PrecioRow -Row 905 -Fecha 45106 -Variedad "Semiduro"   -Calidad "Primera"  -Volumen 600  -PrecioMin 6000  -PrecioMax 7000  -PrecioProm 6500  -PrecioKg 361
Set-PrecioRow -Row 906 -Fecha 45106 -Variedad "Semiduro"   -Calidad "Segunda"  -Volumen 400  -PrecioMin 4000  -PrecioMax 5000  -PrecioProm 4500  -PrecioKg 250

# Ensure the date column keeps the same date-time number format used elsewhere in column D.
$ws.Range("D902:D906").NumberFormat = $ws.Range("D907").NumberFormat
